$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each target cell is written as text (apostrophe-prefixed to stop Excel
# from auto-coercing numeric-looking strings into numbers/dates), then the
# style is reset to "Normal" so the quote-prefix formatting flag does not
# linger as a visible style change on the cell.

$ws.Range('D2').Value = "'51.690.48"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +4.63%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'2.771.06"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +5.66%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  -0.03%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'116.86"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +4.37%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'333.87"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +3.18%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.541"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  +2.74%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = "'  -0.06%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.577"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  +6.54%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'42.06"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +6.27%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'0.0865"
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').Value = "'20.32"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  +3.02%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('E14').Value = "'  +5.42%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'3.206.21"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +5.51%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'2.786.44"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +6.15%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'0.889"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +4.06%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'51.654.28"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +4.72%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'3.33"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +13.99%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'13.53"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  +5.45%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'6.87"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  +3.19%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'0.0₃0978"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  +3.54%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'278.73"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +3.73%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'69.88"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +1.42%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'2.69"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  +6.20%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'26.86"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  +2.63%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'0.999"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  -0.05%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'10.19"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -0.50%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'  +0.53%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = "'  +3.13%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'35.22"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  +1.84%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'50.23"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +1.73%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'5.58"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  +1.91%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'0.0821"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  +1.06%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('B35').Value = "'Celestia"
$ws.Range('B35').Style = 'Normal'
$ws.Range('C35').Value = "'https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range('C35').Style = 'Normal'
$ws.Range('D35').Value = "'19.27"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  +2.27%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('B36').Value = "'FirstDigitalUSD"
$ws.Range('B36').Style = 'Normal'
$ws.Range('C36').Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range('C36').Style = 'Normal'
$ws.Range('D36').Value = "'1.00"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -0.14%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('B37').Value = "'RenderToken"
$ws.Range('B37').Style = 'Normal'
$ws.Range('C37').Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range('C37').Style = 'Normal'
$ws.Range('D37').Value = "'5.04"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  +3.17%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = "'  +3.06%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'3.24"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  +4.95%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = "'  +9.54%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'128.48"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  +0.45%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'23.44"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +5.58%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('B43').Value = "'WEMIXToken"
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').Value = "'2.32"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +7.80%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('B44').Value = "'Stellar"
$ws.Range('B44').Style = 'Normal'
$ws.Range('C44').Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range('C44').Style = 'Normal'
$ws.Range('D44').Value = "'0.115"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +3.65%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('E45').Value = "'  +17.77%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'2.091.59"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +2.21%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'3.31"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  +3.84%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = "'  +3.98%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'5.55"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  +7.19%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = "'  +3.13%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'8.85"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  -0.51%  "
$ws.Range('E51').Style = 'Normal'
